# Convert disabled selection fields to read only text fields.
#
# 1. Insert a new bullet (numPr ilvl=1/numId=3, same list as its
#    neighbours) right after the SSL-in-production bullet (i.e. right
#    before the blank "ind left=1440" spacer paragraph that follows it),
#    describing the newly found issue with clearing a selection field
#    while editing a record.
# 2. Remove the "In MS Explorer disabled and read-only field look
#    differently" bullet entirely (the underlying bug is fixed now that
#    disabled fields are converted to read-only text fields).

$d = $word.ActiveDocument

# --- 1. find the SSL-in-production bullet; the new bullet goes right
#        after it (i.e. right before the following blank paragraph) -------
$findRange = $d.Content
$found = $findRange.Find.Execute(
    "railscast episode 357)",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $sslParaIndex = $findRange.Paragraphs.Item(1).Index
    $targetIndex = $sslParaIndex + 1

    $target = $d.Paragraphs.Item($targetIndex)
    $target.Range.InsertParagraphBefore()

    # re-fetch by index: the freshly inserted (still empty) paragraph now
    # occupies the slot the blank paragraph used to be at
    $newPara = $d.Paragraphs.Item($targetIndex)

    $newParaXml = "<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>" + `
        "<w:pPr>" + `
            "<w:pStyle w:val='ListParagraph'/>" + `
            "<w:numPr><w:ilvl w:val='1'/><w:numId w:val='3'/></w:numPr>" + `
            "<w:jc w:val='both'/>" + `
        "</w:pPr>" + `
        "<w:r><w:t xml:space='preserve'>It is not possible to clear selection field during editing the record, for example </w:t></w:r>" + `
        "<w:proofErr w:type='spellStart'/>" + `
        "<w:r><w:t>accountant_id</w:t></w:r>" + `
        "<w:proofErr w:type='spellEnd'/>" + `
        "<w:r><w:t xml:space='preserve'> field in company table.</w:t></w:r>" + `
    "</w:p>"

    $newPara.Range.InsertXML($newParaXml)
}

# --- 2. delete the "In MS Explorer ..." bullet entirely ---------------------
$findRange2 = $d.Content
$found2 = $findRange2.Find.Execute(
    "In MS Explorer disabled and read-only field look differently",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found2) {
    $paraIndex = $findRange2.Paragraphs.Item(1).Index
    $d.Paragraphs.Item($paraIndex).Range.Delete()
}
